$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "T-shirt"
$ws.Range("B6").Value = "Hassan "
$ws.Range("C6").Value = "Baraka"
